$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.784.62"
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").Value = "1.814.89"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'230.45"
$ws.Range("E5").Value = "  -1.44%  "

$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'39.79"
$ws.Range("E8").Value = "  -9.38%  "

$ws.Range("D9").Value = "'0.322"

$ws.Range("D10").Value = "'0.0681"
$ws.Range("E10").Value = "  -1.81%  "

$ws.Range("E11").Value = "  -2.10%  "

$ws.Range("D12").Value = "2.073.95"
$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").Value = "'11.22"
$ws.Range("E13").Value = "  -1.16%  "

$ws.Range("D14").Value = "'0.666"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").Value = "1.779.82"
$ws.Range("E15").Value = "  -3.21%  "

$ws.Range("D16").Value = "'4.59"
$ws.Range("E16").Value = "  -3.65%  "

$ws.Range("D17").Value = "34.748.90"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").Value = "'69.39"
$ws.Range("E18").Value = "  -1.83%  "

$ws.Range("D19").Value = "0.0₃0780"
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").Value = "'239.83"
$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("D21").Value = "'11.90"
$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("D22").Value = "'4.64"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "'2.24"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").Value = "'173.58"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("D26").Value = "'7.74"
$ws.Range("E26").Value = "  -1.97%  "

$ws.Range("E27").Value = "  +1.38%  "

$ws.Range("D28").Value = "'17.30"
$ws.Range("E28").Value = "  -2.19%  "

$ws.Range("E29").Value = "  -3.98%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").Value = "'3.99"
$ws.Range("E31").Value = "  +0.95%  "

$ws.Range("D32").Value = "'0.0545"
$ws.Range("E32").Value = "  -2.59%  "

$ws.Range("D33").Value = "'3.92"
$ws.Range("E33").Value = "  -3.56%  "

$ws.Range("D34").Value = "'1.25"
$ws.Range("E34").Value = "  +11.58%  "

$ws.Range("D35").Value = "'1.80"
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("D36").Value = "'0.688"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "'91.05"
$ws.Range("E37").Value = "  -4.34%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.35"
$ws.Range("E38").Value = "  +6.15%  "

$ws.Range("D39").Value = "1.330.08"
$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("D41").Value = "'0.971"
$ws.Range("E41").Value = "  -3.72%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.27"
$ws.Range("E42").Value = "  -5.65%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.42"
$ws.Range("E43").Value = "  -1.35%  "

$ws.Range("D44").Value = "'14.21"
$ws.Range("E44").Value = "  -7.34%  "

$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  -1.87%  "

$ws.Range("D46").Value = "'0.0520"
$ws.Range("E46").Value = "  +1.69%  "

$ws.Range("D47").Value = "'6.12"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("D48").Value = "1.992.78"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").Value = "'0.0661"
$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").Value = "'96.99"
$ws.Range("E51").Value = "  -4.24%  "
